$wb = $excel.ActiveWorkbook

# Add the new empty-ER "PRIDE_PROTEOMICS" sheet as the last sheet in the
# workbook (it becomes the active tab, matching the commit that fills an
# empty ER sheet into every template except Imaging).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PRIDE_PROTEOMICS"

# Content mirrors the main annotation table's columns (Source Name,
# Parameter [acquisition software], Parameter [analysis software],
# Parameter [data processing software], Data File Name) together with the
# TermSourceRef/Ontology/TAN scaffold columns used by the ER sheets.
$data = @(
    @("", "TermSourceRef", "Ontology", "TAN", "Content type (validation)", "Notes during templating", "Target term", "Instruction", "Requirement (m/o/n)", "Value (cv/s/d)", "Additional information", "Review comments"),
    @("Source Name", "", "", "", "", "", "", "", "", "", "", ""),
    @("Parameter [acquisition software]", "MS:1001455", "MS", "http://purl.obolibrary.org/obo/MS_1001455", "", "", "", "", "", "", "", ""),
    @("Parameter [analysis software]", "MS:1001456", "MS", "http://purl.obolibrary.org/obo/MS_1001456", "", "", "", "", "", "", "", ""),
    @("Parameter [data processing software]", "MS:1001457", "MS", "http://purl.obolibrary.org/obo/MS_1001457", "", "", "", "", "", "", "", ""),
    @("Data File Name", "", "", "", "", "", "", "", "", "", "", "")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($value -ne "") {
            $newSheet.Cells.Item($r + 1, $c + 1).Value = $value
        }
    }
}

# Column widths matching the bestFit widths of the source sheet.
$widths = @(34.7109375, 14.5703125, 9.140625, 40.28515625, 23.5703125, 22.85546875, 11.28515625, 10.5703125, 20.42578125, 13.5703125, 21.5703125, 17.42578125)
for ($c = 0; $c -lt $widths.Length; $c++) {
    $newSheet.Columns.Item($c + 1).ColumnWidth = $widths[$c]
}

# Select the whole sheet (matches the authored sqref="A1:XFD1048576") and
# make this new sheet the active tab, so it is the one shown on open.
$newSheet.Cells.Select() | Out-Null
$newSheet.Activate()
